$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Data updates on the "Training Dashboard" sheet (rows 3 & 4)
#    H3/H4: -19610 -> -19618
#    I3/I4: "08-Sep-2025" -> "16-Sep-2025" (kept as literal text, not a
#    parsed date, matching the rest of the "wrongly formatted date" data)
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("H3").Value = -19618
$ws1.Range("H4").Value = -19618

$ws1.Range("I3").Value = "'16-Sep-2025"
$ws1.Range("I4").Value = "'16-Sep-2025"

# ----------------------------------------------------------------------
# 2) Style updates: the bold title + the bold header row now render in
#    white instead of the default (black) font color. The title also
#    drops its old 14pt override, falling back to the standard 11pt
#    size used throughout the rest of the workbook.
# ----------------------------------------------------------------------
$white = 16777215

$ws1Title = $ws1.Range("A1")
$ws1Title.Font.Bold = $true
$ws1Title.Font.Size = 11
$ws1Title.Font.Color = $white

$ws1Header = $ws1.Range("A2:K2")
$ws1Header.Font.Bold = $true
$ws1Header.Font.Color = $white

$ws2 = $wb.Worksheets.Item(2)

$ws2Title = $ws2.Range("A1")
$ws2Title.Font.Bold = $true
$ws2Title.Font.Size = 11
$ws2Title.Font.Color = $white

$ws2Header = $ws2.Range("A2:G2")
$ws2Header.Font.Bold = $true
$ws2Header.Font.Color = $white
